# Update workbook "广州-漫展信息.xlsx" with freshly scraped "想去人数" (interest
# count) figures, and remove a duplicated row from the "全部类型" (all types)
# summary sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - update F column (想去人数) values
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 53
$ws1.Range("F3").Value = 275
$ws1.Range("F4").Value = 970
$ws1.Range("F6").Value = 437
$ws1.Range("F7").Value = 659
$ws1.Range("F11").Value = 378
$ws1.Range("F12").Value = 178
$ws1.Range("F13").Value = 36
$ws1.Range("F14").Value = 764
$ws1.Range("F16").Value = 1903
$ws1.Range("F17").Value = 418
$ws1.Range("F18").Value = 5671
$ws1.Range("F19").Value = 416
$ws1.Range("F20").Value = 507
$ws1.Range("F22").Value = 75
$ws1.Range("F23").Value = 7
$ws1.Range("F24").Value = 181

# ------------------------------------------------------------------
# Sheet "演出" (Performances) - update F column (想去人数) values
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 26
$ws2.Range("F7").Value = 495
$ws2.Range("F13").Value = 109
$ws2.Range("F15").Value = 7

# ------------------------------------------------------------------
# Sheet "本地生活" (Local life) - update F column (想去人数) values
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5420
$ws3.Range("F3").Value = 358
$ws3.Range("F4").Value = 343

# ------------------------------------------------------------------
# Sheet "全部类型" (All types) - this sheet aggregates all rows from the
# other sheets. Row 13 ("广州·神山羊2024巡演ENCOUNTER") was an exact
# duplicate of row 12, so it is removed; all following rows shift up by
# one (row 49 disappears, dimension becomes A1:I48). Afterwards the same
# F column (想去人数) values are refreshed at their new row positions.
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows("13:13").Delete()

$ws4.Range("F2").Value = 53
$ws4.Range("F3").Value = 5420
$ws4.Range("F4").Value = 358
$ws4.Range("F6").Value = 343
$ws4.Range("F7").Value = 275
$ws4.Range("F10").Value = 26
$ws4.Range("F12").Value = 495
$ws4.Range("F13").Value = 970
$ws4.Range("F17").Value = 437
$ws4.Range("F18").Value = 659
$ws4.Range("F23").Value = 378
$ws4.Range("F24").Value = 178
$ws4.Range("F26").Value = 36
$ws4.Range("F28").Value = 764
$ws4.Range("F30").Value = 109
$ws4.Range("F31").Value = 1903
$ws4.Range("F32").Value = 418
$ws4.Range("F33").Value = 5671
$ws4.Range("F35").Value = 416
$ws4.Range("F36").Value = 507
$ws4.Range("F38").Value = 75
$ws4.Range("F39").Value = 7
$ws4.Range("F40").Value = 7
$ws4.Range("F41").Value = 181

Write-Output "Applied visitor count updates across sheets and removed duplicate row from the all-types summary sheet."
